$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNumber, new DAMSLTag (col I), new DialogAct (col J)
$updates = @(
    @{ Row = 32;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 34;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 49;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 67;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 68;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 69;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 93;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 102; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 110; Tag = "sd"; Act = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Range("I$($u.Row)").Value = $u.Tag
    $ws.Range("J$($u.Row)").Value = $u.Act
}
